$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the added Q8 column - copy the same header style (bold,
# bordered, centered) as the neighboring header cells before setting its value.
$ws.Cells.Item(1, 9).Copy()
$ws.Cells.Item(1, 10).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(1, 10).Value = "Q8"

# Row 4 gains values for columns G:J
$ws.Cells.Item(4, 7).Value = 0.3837439847994517
$ws.Cells.Item(4, 8).Value = -1.436992267669069
$ws.Cells.Item(4, 9).Value = 0.3765140280931547
$ws.Cells.Item(4, 10).Value = -0.2577465226711695

# Row 8 gains values for columns G:I
$ws.Cells.Item(8, 7).Value = 1.061234703579416
$ws.Cells.Item(8, 8).Value = 0.6797930702803257
$ws.Cells.Item(8, 9).Value = 0.3873308289134342
